# Apply the Jan 2 2023 symbol-list update to Sheet1 (cryptos.xlsx export).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $value) {
    # Force literal text (even for number-/percent-looking strings) the same way
    # a user typing a leading apostrophe in Excel would, then drop the resulting
    # quote-prefix style so the cell format matches its original (General) style.
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $value
    $cell.ClearFormats()
}

Set-TextCell "D2" "246.53"
Set-TextCell "E2" "0.61%"
Set-TextCell "D3" "29.83"
Set-TextCell "E3" "9.69%"
Set-TextCell "D4" "5.179"
Set-TextCell "E4" "1.67%"
Set-TextCell "D5" "0.05731"
Set-TextCell "E5" "0.46%"
Set-TextCell "E6" "1.18%"
Set-TextCell "D7" "0.8565"
Set-TextCell "E7" "4.46%"
Set-TextCell "D8" "0.8665"
Set-TextCell "E8" "0.61%"
Set-TextCell "E9" "2.36%"
Set-TextCell "E10" "2.08%"
Set-TextCell "D11" "0.02944"
Set-TextCell "E11" "4.20%"
Set-TextCell "D12" "0.09390"
Set-TextCell "E12" "-0.15%"
Set-TextCell "D13" "0.001525"
Set-TextCell "E13" "0.88%"
Set-TextCell "D14" "0.04178"
Set-TextCell "E14" "3.03%"
Set-TextCell "D15" "0.0005979"
Set-TextCell "E15" "-0.62%"
Set-TextCell "D16" "0.005996"
Set-TextCell "E16" "-2.43%"
Set-TextCell "E17" "5,071.82%"
Set-TextCell "D18" "3.487"
Set-TextCell "E18" "-0.57%"
Set-TextCell "D19" "3.099"
Set-TextCell "E19" "2.98%"
Set-TextCell "D20" "2.189"
Set-TextCell "E20" "-1.80%"
Set-TextCell "D22" "0.03443"
Set-TextCell "E22" "7.70%"
Set-TextCell "E23" "1.07%"
Set-TextCell "D24" "3.486"
Set-TextCell "E24" "-1.93%"
Set-TextCell "E25" "0.47%"
Set-TextCell "D26" "0.005019"
Set-TextCell "E26" "12.02%"
Set-TextCell "D27" "0.001228"
Set-TextCell "E27" "1.08%"
Set-TextCell "E40" "0.77%"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextCell "D41" "0.1074"
Set-TextCell "E41" "1.38%"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextCell "D42" "0.002539"
Set-TextCell "E42" "2.05%"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextCell "D43" "0.003489"
Set-TextCell "E43" "-41.13%"
Set-TextCell "D44" "0.009664"
Set-TextCell "E44" "2.99%"
Set-TextCell "D45" "0.00005221"
Set-TextCell "E45" "1.40%"
Set-TextCell "E46" "0.05%"
Set-TextCell "E47" "-45.05%"
Set-TextCell "D48" "0.002520"
Set-TextCell "E48" "0.03%"
Set-TextCell "D49" "0.00002100"
Set-TextCell "E49" "0.05%"
Set-TextCell "E50" "0.05%"
